$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2364.2222
$ws.Range("I6").Value = 79.57143000000001
$ws.Range("J6").Value = 3818.0908
$ws.Range("K6").Value = 238.71429
$ws.Range("L6").Value = 11454.2724
$ws.Range("M6").Value = -126.71429
$ws.Range("N6").Value = -11678.2724
$ws.Range("H9").Value = 199.37038
$ws.Range("I9").Value = 199.37038
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 199.37038
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -30.37038000000001
$ws.Range("N9").ClearContents()
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H105").Value = 39715.832
$ws.Range("J105").Value = 39715.832
$ws.Range("L105").Value = 39715.832
$ws.Range("N105").Value = -46703.832

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3625.9092
$ws.Range("I122").Value = 1005
$ws.Range("K122").Value = 3015
$ws.Range("M122").Value = -565

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 441.66666
$ws.Range("I7").Value = 560
$ws.Range("J7").Value = 293.75
$ws.Range("K7").Value = 560
$ws.Range("L7").Value = 293.75
$ws.Range("M7").Value = -447
$ws.Range("N7").Value = -519.75
$ws.Range("H23").Value = 35000
$ws.Range("J23").Value = 35000
$ws.Range("L23").Value = 35000
$ws.Range("N23").Value = -35480
$ws.Range("H27").Value = 35000
$ws.Range("J27").Value = 35000
$ws.Range("L27").Value = 35000
$ws.Range("N27").Value = -35384
$ws.Range("H58").Value = 2889.6865
$ws.Range("I58").Value = 1771.9
$ws.Range("J58").Value = 6177.294
$ws.Range("K58").Value = 1771.9
$ws.Range("L58").Value = 6177.294
$ws.Range("M58").Value = -1568.9
$ws.Range("N58").Value = -6583.294
$ws.Range("H105").Value = 1797.2
$ws.Range("I105").Value = 1119.75
$ws.Range("J105").Value = 2571.4285
$ws.Range("K105").Value = 1119.75
$ws.Range("L105").Value = 2571.4285
$ws.Range("M105").Value = 627.25
$ws.Range("N105").Value = -6065.4285
$ws.Range("H136").Value = 2889.6865
$ws.Range("I136").Value = 1771.9
$ws.Range("J136").Value = 6177.294
$ws.Range("K136").Value = 5315.700000000001
$ws.Range("L136").Value = 18531.882
$ws.Range("M136").Value = -2765.700000000001
$ws.Range("N136").Value = -23631.882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1478
$ws.Range("J22").Value = 1478
$ws.Range("L22").Value = 4434
$ws.Range("N22").Value = -4772
$ws.Range("H27").Value = 1478
$ws.Range("J27").Value = 1478
$ws.Range("L27").Value = 4434
$ws.Range("N27").Value = -4638
$ws.Range("H38").Value = 171.71428
$ws.Range("J38").Value = 185.71428
$ws.Range("L38").Value = 557.14284
$ws.Range("N38").Value = -1251.14284
$ws.Range("H113").Value = 2551675.5
$ws.Range("I113").Value = 564.54285
$ws.Range("J113").Value = 8929453
$ws.Range("K113").Value = 1693.62855
$ws.Range("L113").Value = 26788359
$ws.Range("M113").Value = 476.3714499999999
$ws.Range("N113").Value = -26792699
$ws.Range("H117").Value = 1579.8334
$ws.Range("J117").Value = 1790
$ws.Range("L117").Value = 5370
$ws.Range("N117").Value = -12254

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H122").Value = 7734.1113
$ws.Range("I122").Value = 2121.6
$ws.Range("J122").Value = 14749.75
$ws.Range("K122").Value = 6364.799999999999
$ws.Range("L122").Value = 44249.25
$ws.Range("M122").Value = -3914.799999999999
$ws.Range("N122").Value = -49149.25
$ws.Range("H126").Value = 3529.293
$ws.Range("I126").Value = 2971.831
$ws.Range("J126").Value = 4942.857
$ws.Range("K126").Value = 8915.493
$ws.Range("L126").Value = 14828.571
$ws.Range("M126").Value = -6445.493
$ws.Range("N126").Value = -19768.571
$ws.Range("H132").Value = 4654.8125
$ws.Range("I132").Value = 3379.9092
$ws.Range("J132").Value = 7459.6
$ws.Range("K132").Value = 10139.7276
$ws.Range("L132").Value = 22378.8
$ws.Range("M132").Value = -7609.7276
$ws.Range("N132").Value = -27438.8
$ws.Range("H135").Value = 24163.637
$ws.Range("J135").Value = 24163.637
$ws.Range("L135").Value = 24163.637
$ws.Range("N135").Value = -34303.637

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6865.5835
$ws.Range("J7").Value = 7533.1665
$ws.Range("L7").Value = 7533.1665
$ws.Range("N7").Value = -7757.1665
$ws.Range("H40").Value = 4699.7144
$ws.Range("I40").Value = 3279.8
$ws.Range("J40").Value = 8249.5
$ws.Range("K40").Value = 3279.8
$ws.Range("L40").Value = 8249.5
$ws.Range("M40").Value = -3143.8
$ws.Range("N40").Value = -8521.5
$ws.Range("H82").Value = 1509.875
$ws.Range("I82").Value = 689.7857
$ws.Range("J82").Value = 2658
$ws.Range("K82").Value = 689.7857
$ws.Range("L82").Value = 2658
$ws.Range("M82").Value = -328.7857
$ws.Range("N82").Value = -3380
$ws.Range("H85").Value = 1509.875
$ws.Range("I85").Value = 689.7857
$ws.Range("J85").Value = 2658
$ws.Range("K85").Value = 689.7857
$ws.Range("L85").Value = 2658
$ws.Range("M85").Value = 558.2143
$ws.Range("N85").Value = -5154
$ws.Range("H126").Value = 6865.5835
$ws.Range("J126").Value = 7533.1665
$ws.Range("L126").Value = 22599.4995
$ws.Range("N126").Value = -27539.4995
$ws.Range("H136").Value = 5108.778
$ws.Range("I136").Value = 1457.25
$ws.Range("J136").Value = 8030
$ws.Range("K136").Value = 4371.75
$ws.Range("L136").Value = 24090
$ws.Range("M136").Value = -1821.75
$ws.Range("N136").Value = -29190
$ws.Range("H140").Value = 65499.625
$ws.Range("J140").Value = 65499.625
$ws.Range("L140").Value = 65499.625
$ws.Range("N140").Value = -75859.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4920.5557
$ws.Range("I126").Value = 2146
$ws.Range("J126").Value = 10469.667
$ws.Range("K126").Value = 6438
$ws.Range("L126").Value = 31409.001
$ws.Range("M126").Value = -3968
$ws.Range("N126").Value = -36349.001
$ws.Range("H136").Value = 6273.1177
$ws.Range("I136").Value = 2604.75
$ws.Range("K136").Value = 7814.25
$ws.Range("M136").Value = -5264.25

Write-Host "Applied scheduled Chocobo Profits market-data refresh to ALC, ARM, CRP, CUL, GSM, LTW, WVR"
